# [FEATURE] Migracion casos de UFT
# Adds 6 new rows (49-54) of test-case data to the "Users" sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Users")

# --- New rows ---------------------------------------------------------
# Row 49: F00019 / 019
$ws.Cells.Item(49, 1).Value = "F00019"
$ws.Cells.Item(49, 3).NumberFormat = "@"
$ws.Cells.Item(49, 3).HorizontalAlignment = -4152
$ws.Cells.Item(49, 3).Value = "019"

# Row 50: F02171 / 019
$ws.Cells.Item(50, 1).Value = "F02171"
$ws.Cells.Item(50, 3).NumberFormat = "@"
$ws.Cells.Item(50, 3).HorizontalAlignment = -4152
$ws.Cells.Item(50, 3).Value = "019"

# Row 51: F00103 / 103 (plain number)
$ws.Cells.Item(51, 1).Value = "F00103"
$ws.Cells.Item(51, 3).Value = 103

# Row 52: F02582 / 103 (plain number)
$ws.Cells.Item(52, 1).Value = "F02582"
$ws.Cells.Item(52, 3).Value = 103

# Row 53: FMASTROIANNI / Casa central
$ws.Cells.Item(53, 1).Value = "FMASTROIANNI"
$ws.Cells.Item(53, 3).NumberFormat = "@"
$ws.Cells.Item(53, 3).HorizontalAlignment = -4152
$ws.Cells.Item(53, 3).Value = "Casa central"

# Row 54: MABRUNI / 102 (plain number)
$ws.Cells.Item(54, 1).Value = "MABRUNI"
$ws.Cells.Item(54, 3).Value = 102

# --- View state: select the new last cell, matching the authored edit.
$ws.Activate()
$ws.Cells.Item(54, 3).Select()
